$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the daily log. Insert a fresh row
# above the existing row 187 (pushing rows 187:264 down to 188:265) and
# populate it with the new observation.
$ws.Rows.Item(187).Insert()

$ws.Cells.Item(187, 1).Value = 4
$ws.Cells.Item(187, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(187, 3).Value = "Los Lagos"
$ws.Cells.Item(187, 4).Value = 44704
$ws.Cells.Item(187, 5).Value = 10
$ws.Cells.Item(187, 6).Value = 100112043
$ws.Cells.Item(187, 7).Value = "Pepino ensalada"
$ws.Cells.Item(187, 8).Value = "Sin especificar"
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 120
$ws.Cells.Item(187, 11).Value = 22000
$ws.Cells.Item(187, 12).Value = 22000
$ws.Cells.Item(187, 13).Value = 22000
$ws.Cells.Item(187, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(187, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(187, 16).Value = 367
$ws.Cells.Item(187, 17).Value = 60
$ws.Cells.Item(187, 18).Value = "Hortaliza"
